$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.140395522117615
$ws.Range("B1").Value = 2.228063821792603
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.183600425720215
$ws.Range("E1").Value = 1.075807690620422
